$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 43
$ws.Range("F4").Value = 261
$ws.Range("F5").Value = 2010
$ws.Range("F6").Value = 81
$ws.Range("F7").Value = 481
$ws.Range("F8").Value = 425
$ws.Range("F9").Value = 219
$ws.Range("F10").Value = 7247
$ws.Range("F11").Value = 206
$ws.Range("F13").Value = 452
$ws.Range("F15").Value = 3118
$ws.Range("F16").Value = 1797
$ws.Range("F17").Value = 159
$ws.Range("F18").Value = 3
$ws.Range("F20").Value = 108
$ws.Range("F21").Value = 172
$ws.Range("F24").Value = 176
$ws.Range("F26").Value = 984
$ws.Range("F27").Value = 197
$ws.Range("F28").Value = 4131

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 77
$ws.Range("F3").Value = 29

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 128
$ws.Range("F3").Value = 716

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 128
$ws.Range("F4").Value = 716
$ws.Range("F5").Value = 43
$ws.Range("F6").Value = 77
$ws.Range("F7").Value = 261
$ws.Range("F8").Value = 2010
$ws.Range("F9").Value = 29
$ws.Range("F11").Value = 81
$ws.Range("F12").Value = 481
$ws.Range("F13").Value = 425
$ws.Range("F14").Value = 219
$ws.Range("F15").Value = 7247
$ws.Range("F16").Value = 206
$ws.Range("F18").Value = 452
$ws.Range("F20").Value = 3118
$ws.Range("F21").Value = 1797
$ws.Range("F22").Value = 159
$ws.Range("F23").Value = 3
$ws.Range("F25").Value = 108
$ws.Range("F26").Value = 172
$ws.Range("F29").Value = 176
$ws.Range("F31").Value = 984
$ws.Range("F32").Value = 197
$ws.Range("F33").Value = 4131
